# The deck originally used the "Integral" theme (colour scheme) for its
# slide master / presentation theme. The authored edit swaps the theme
# applied to the deck so it instead uses the classic "Office Theme"
# colour scheme (while the "Integral" colours move to become the
# secondary/unused theme part).
#
# The PowerPoint object model exposes the applied theme's 12 theme
# colours through Theme.ThemeColorScheme (Dark1, Light1, Dark2, Light2,
# Accent1-6, Hyperlink, FollowedHyperlink, in that order). Re-pointing
# every slot to the "Office Theme" palette reproduces the colour swap.

function HexToComRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (the deck's original/default theme palette).
$officeColors = @(
    "000000",  # 1  Dark1
    "FFFFFF",  # 2  Light1
    "44546A",  # 3  Dark2
    "E7E6E6",  # 4  Light2
    "5B9BD5",  # 5  Accent1
    "ED7D31",  # 6  Accent2
    "A5A5A5",  # 7  Accent3
    "FFC000",  # 8  Accent4
    "4472C4",  # 9  Accent5
    "70AD47",  # 10 Accent6
    "0563C1",  # 11 Hyperlink
    "954F72"   # 12 FollowedHyperlink
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToComRGB $officeColors[$i - 1]
}
